$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.0985
$ws.Range("G2").Value = 0.01789158950617284
$ws.Range("H2").Value = 0.01789158950617284
$ws.Range("I2").Value = -0.1508646688685669
$ws.Range("J2").Value = -0.1508646688685669
$ws.Range("K2").Value = -68.67100000000001
$ws.Range("L2").Value = -0.1655840084876543
$ws.Range("M2").Value = 12.7
$ws.Range("N2").Value = 0.04740574841358716
$ws.Range("O2").Value = -0.1849397853533515
$ws.Range("S2").Value = 12.7
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 11.6
$ws.Range("V2").Value = 0.04329973870847331
$ws.Range("W2").Value = -0.07834865692051016
$ws.Range("X2").Value = 0.08909424715512393
$ws.Range("Y2").Value = -0.1674429040756341
$ws.Range("Z2").Value = 0.6998100223746121
$ws.Range("AA2").Value = -0.0522852562940653
$ws.Range("AB2").Value = 0.07820086300925749
$ws.Range("AC2").Value = -0.1304861193033228
$ws.Range("AD2").Value = 94.416
$ws.Range("AE2").Value = 0.9179773658602803
$ws.Range("AF2").Value = 95.33397736586029
$ws.Range("AG2").Value = 83.7339773658603
$ws.Range("AH2").Value = 0.2624588648265055
$ws.Range("AI2").Value = 0.1798503254953416
$ws.Range("AJ2").Value = 0.238128232069959
$ws.Range("AK2").Value = 0.1615008293979884
$ws.Range("AL2").Value = 4.29
$ws.Range("AM2").Value = 4.29
$ws.Range("AN2").Value = -1.513731903227358
$ws.Range("AO2").Value = -14.70815850815851
$ws.Range("AP2").Value = -1.342471540023092
$ws.Range("AQ2").Value = -14.70815850815851
$ws.Range("B3").Value = "Oxbridge Re Holdings Limited (NasdaqCM:OXBR)"
$ws.Range("D3").Value = -0.319
$ws.Range("G3").Value = -0.125
$ws.Range("H3").Value = -0.125
$ws.Range("I3").Value = 0.00919175721831896
$ws.Range("J3").Value = 0.00919175721831896
$ws.Range("K3").Value = -0.171
$ws.Range("L3").Value = -0.1526785714285714
$ws.Range("U3").Value = 3.44
$ws.Range("V3").Value = 0.3214953271028038
$ws.Range("W3").Value = -0.02145545796737767
$ws.Range("X3").Value = 0.07792933076350653
$ws.Range("Y3").Value = -0.09938478873088419
$ws.Range("Z3").Value = 0.2138005931214724
$ws.Range("AA3").Value = 0.001965203145105169
$ws.Range("AB3").Value = 0.07656287168585088
$ws.Range("AC3").Value = -0.07459766854074572
$ws.Range("AD3").Value = 0.216
$ws.Range("AE3").Value = 0.1385261595774138
$ws.Range("AF3").Value = 0.3545261595774138
$ws.Range("AG3").Value = -3.085473840422586
$ws.Range("AH3").Value = 0.03207067896530867
$ws.Range("AI3").Value = 0.04326377787726735
$ws.Range("AJ3").Value = -0.4052089093609236
$ws.Range("AK3").Value = -0.6489550665753042
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AN3").Value = 4.5
$ws.Range("AO3").ClearContents()
$ws.Range("AP3").Value = -64.28070500880388
$ws.Range("AQ3").ClearContents()
$ws.Range("B4").Value = "Greenlight Capital Re, Ltd. (NasdaqGS:GLRE)"
$ws.Range("D4").Value = 0.122
$ws.Range("G4").Value = 0.01827852998065764
$ws.Range("H4").Value = 0.01827852998065764
$ws.Range("I4").Value = -0.1512980905252818
$ws.Range("J4").Value = -0.1512980905252818
$ws.Range("K4").Value = -68.5
$ws.Range("L4").Value = -0.1656189555125725
$ws.Range("M4").Value = 12.7
$ws.Range("N4").Value = 0.04937791601866252
$ws.Range("O4").Value = -0.1854014598540146
$ws.Range("S4").Value = 12.7
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 8.16
$ws.Range("V4").Value = 0.03172628304821151
$ws.Range("W4").Value = -0.1352418558736427
$ws.Range("X4").Value = 0.1002591635467413
$ws.Range("Y4").Value = -0.235501019420384
$ws.Range("Z4").Value = 0.7041444830094117
$ws.Range("AA4").Value = -0.1065357157332358
$ws.Range("AB4").Value = 0.07983885433266409
$ws.Range("AC4").Value = -0.1863745700658999
$ws.Range("AD4").Value = 94.2
$ws.Range("AE4").Value = 0.7794512062828665
$ws.Range("AF4").Value = 94.97945120628287
$ws.Range("AG4").Value = 86.81945120628288
$ws.Range("AH4").Value = 0.2696904969354681
$ws.Range("AI4").Value = 0.1819950009274085
$ws.Range("AJ4").Value = 0.2523678556600676
$ws.Range("AK4").Value = 0.1690016817592151
$ws.Range("AL4").Value = 4.29
$ws.Range("AM4").Value = 4.29
$ws.Range("AN4").Value = -1.509107511895035
$ws.Range("AO4").Value = -14.70862470862471
$ws.Range("AP4").Value = -1.390869278068004
$ws.Range("AQ4").Value = -14.70862470862471
